$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45205 -> 45206, i.e. 2023-10-06 -> 2023-10-07) for every data row
# (rows 2 through 236).
for ($row = 2; $row -le 236; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
